$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("JSS 3A")

# Fix C4: was stored as text "38", should be a real number 38
$ws.Range("C4").Value = 38

# Append the new submission row (row 5)
$ws.Range("A5").Value = "2026-02-08 18:57:43"
$ws.Range("B5").Value = "Halima Sadiya Abubakar"

# Admission No for this submission needs to stay text (not auto-coerced
# to a number like C4 was before this sync), so force text formatting,
# write it, then drop back to the default style (keeps the text type
# without leaving a visible number-format override on the cell).
$ws.Range("C5").NumberFormat = "@"
$ws.Range("C5").Value = "25"
$ws.Range("C5").Style = "Normal"

$ws.Range("D5").Value = 9
